$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.334522
$ws.Range("H2").Value = 4.003566
$ws.Range("I2").Value = 0.4120903366177529
$ws.Range("J2").Value = 0.4120903366177529
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06447966666666667
$ws.Range("N2").Value = 0.193439
$ws.Range("O2").Value = 0.001101138907643723
$ws.Range("P2").Value = 0.001101138907643722
$ws.Range("Q2").Value = 0.08604953371933334
$ws.Range("R2").Value = 0.774445803474
$ws.Range("S2").Value = 0.0004537687031138063
$ws.Range("T2").Value = 0.0004537687031138063

# Row 3
$ws.Range("G3").Value = 1.334522
$ws.Range("H3").Value = 4.003566
$ws.Range("I3").Value = 0.4120903366177529
$ws.Range("J3").Value = 0.4120903366177529
$ws.Range("O3").Value = 0.00657695954769643
$ws.Range("P3").Value = 0.006576959547696431
$ws.Range("Q3").Value = 0.5139626784973333
$ws.Range("R3").Value = 4.625664106476
$ws.Range("S3").Value = 0.002710301473931566
$ws.Range("T3").Value = 0.002710301473931566

# Row 4
$ws.Range("G4").Value = 1.334522
$ws.Range("H4").Value = 4.003566
$ws.Range("I4").Value = 0.4120903366177529
$ws.Range("J4").Value = 0.4120903366177529
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008175
$ws.Range("N4").Value = 0.024525
$ws.Range("O4").Value = 0.0001396069650378791
$ws.Range("P4").Value = 0.0001396069650378791
$ws.Range("Q4").Value = 0.01090971735
$ws.Range("R4").Value = 0.09818745615
$ws.Range("S4").Value = 0.00005753068121664245
$ws.Range("T4").Value = 0.00005753068121664245

# Row 5
$ws.Range("G5").Value = 1.334522
$ws.Range("H5").Value = 4.003566
$ws.Range("I5").Value = 0.4120903366177529
$ws.Range("J5").Value = 0.4120903366177529
$ws.Range("M5").Value = 58.099467
$ws.Range("N5").Value = 174.298401
$ws.Range("O5").Value = 0.992182294579622
$ws.Range("P5").Value = 0.992182294579622
$ws.Range("Q5").Value = 77.53501689977401
$ws.Range("R5").Value = 697.8151520979661
$ws.Range("S5").Value = 0.4088687357594909
$ws.Range("T5").Value = 0.4088687357594909

# Row 6
$ws.Range("I6").Value = 0.4618070555578372
$ws.Range("J6").Value = 0.4618070555578372
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.06447966666666667
$ws.Range("N6").Value = 0.193439
$ws.Range("O6").Value = 0.001101138907643723
$ws.Range("P6").Value = 0.001101138907643722
$ws.Range("Q6").Value = 0.09643099647811113
$ws.Range("R6").Value = 0.867878968303
$ws.Range("S6").Value = 0.0005085137166991208
$ws.Range("T6").Value = 0.0005085137166991207

# Row 7
$ws.Range("I7").Value = 0.4618070555578372
$ws.Range("J7").Value = 0.4618070555578372
$ws.Range("O7").Value = 0.00657695954769643
$ws.Range("P7").Value = 0.006576959547696431
$ws.Range("S7").Value = 0.003037286323244693
$ws.Range("T7").Value = 0.003037286323244694

# Row 8
$ws.Range("I8").Value = 0.4618070555578372
$ws.Range("J8").Value = 0.4618070555578372
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.008175
$ws.Range("N8").Value = 0.024525
$ws.Range("O8").Value = 0.0001396069650378791
$ws.Range("P8").Value = 0.0001396069650378791
$ws.Range("Q8").Value = 0.012225922325
$ws.Range("R8").Value = 0.110033300925
$ws.Range("S8").Value = 0.00006447148145950886
$ws.Range("T8").Value = 0.00006447148145950887

# Row 9
$ws.Range("I9").Value = 0.4618070555578372
$ws.Range("J9").Value = 0.4618070555578372
$ws.Range("M9").Value = 58.099467
$ws.Range("N9").Value = 174.298401
$ws.Range("O9").Value = 0.992182294579622
$ws.Range("P9").Value = 0.992182294579622
$ws.Range("Q9").Value = 86.88924411815302
$ws.Range("R9").Value = 782.0031970633771
$ws.Range("S9").Value = 0.4581967840364339
$ws.Range("T9").Value = 0.4581967840364339

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4083733333333333
$ws.Range("H10").Value = 1.22512
$ws.Range("I10").Value = 0.1261026078244099
$ws.Range("J10").Value = 0.1261026078244099
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.06447966666666667
$ws.Range("N10").Value = 0.193439
$ws.Range("O10").Value = 0.001101138907643723
$ws.Range("P10").Value = 0.001101138907643722
$ws.Range("Q10").Value = 0.02633177640888889
$ws.Range("R10").Value = 0.23698598768
$ws.Range("S10").Value = 0.0001388564878307954
$ws.Range("T10").Value = 0.0001388564878307954

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.4083733333333333
$ws.Range("H11").Value = 1.22512
$ws.Range("I11").Value = 0.1261026078244099
$ws.Range("J11").Value = 0.1261026078244099
$ws.Range("O11").Value = 0.00657695954769643
$ws.Range("P11").Value = 0.006576959547696431
$ws.Range("Q11").Value = 0.1572762773688889
$ws.Range("R11").Value = 1.41548649632
$ws.Range("S11").Value = 0.0008293717505201712
$ws.Range("T11").Value = 0.0008293717505201714

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.4083733333333333
$ws.Range("H12").Value = 1.22512
$ws.Range("I12").Value = 0.1261026078244099
$ws.Range("J12").Value = 0.1261026078244099
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.008175
$ws.Range("N12").Value = 0.024525
$ws.Range("O12").Value = 0.0001396069650378791
$ws.Range("P12").Value = 0.0001396069650378791
$ws.Range("Q12").Value = 0.003338452
$ws.Range("R12").Value = 0.030046068
$ws.Range("S12").Value = 0.00001760480236172777
$ws.Range("T12").Value = 0.00001760480236172777

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.4083733333333333
$ws.Range("H13").Value = 1.22512
$ws.Range("I13").Value = 0.1261026078244099
$ws.Range("J13").Value = 0.1261026078244099
$ws.Range("M13").Value = 58.099467
$ws.Range("N13").Value = 174.298401
$ws.Range("O13").Value = 0.992182294579622
$ws.Range("P13").Value = 0.992182294579622
$ws.Range("Q13").Value = 23.72627300368
$ws.Range("R13").Value = 213.53645703312
$ws.Range("S13").Value = 0.1251167747836972
$ws.Range("T13").Value = 0.1251167747836972

